# Auto-generated: applies the "Updated cryptos list" diff to Sheet1.
# Only D-column values that parse as plain numbers get a leading
# apostrophe so Excel keeps them as literal text (matching the source
# workbook, where these cells are inline strings, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.973.93"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "2.752.60"
$ws.Range("E3").Value = "  +0.10%  "

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").Value = "'576.86"
$ws.Range("E5").Value = "  -0.66%  "

$ws.Range("D6").Value = "'157.58"
$ws.Range("E6").Value = "  +2.73%  "

$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("E8").Value = "  -0.25%  "

$ws.Range("D9").Value = "'0.110"
$ws.Range("E9").Value = "  -1.93%  "

$ws.Range("D10").Value = "'5.82"
$ws.Range("E10").Value = "  -13.63%  "

$ws.Range("D11").Value = "'0.385"
$ws.Range("E11").Value = "  -1.22%  "

$ws.Range("D12").Value = "'0.157"
$ws.Range("E12").Value = "  -2.09%  "

$ws.Range("D13").Value = "3.240.59"
$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("D14").Value = "'26.75"
$ws.Range("E14").Value = "  +1.76%  "

$ws.Range("D15").Value = "63.910.14"
$ws.Range("E15").Value = "  +0.38%  "

$ws.Range("D16").Value = "'0.0000152"
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").Value = "2.756.61"
$ws.Range("E17").Value = "  +0.46%  "

$ws.Range("D18").Value = "'12.06"
$ws.Range("E18").Value = "  +0.70%  "

$ws.Range("D19").Value = "'4.87"
$ws.Range("E19").Value = "  -0.06%  "

$ws.Range("D20").Value = "'359.14"
$ws.Range("E20").Value = "  -0.33%  "

$ws.Range("D21").Value = "'6.78"
$ws.Range("E21").Value = "  -3.03%  "

$ws.Range("D22").Value = "'0.551"
$ws.Range("E22").Value = "  +2.38%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").Value = "'66.14"
$ws.Range("E24").Value = "  +0.40%  "

$ws.Range("D25").Value = "'0.170"
$ws.Range("E25").Value = "  +1.37%  "

$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'8.44"
$ws.Range("E27").Value = "  -1.54%  "

$ws.Range("D28").Value = "0.0₃0925"
$ws.Range("E28").Value = "  +3.46%  "

$ws.Range("D29").Value = "'1.95"
$ws.Range("E29").Value = "  -2.70%  "

$ws.Range("D30").Value = "'7.00"
$ws.Range("E30").Value = "  -1.16%  "

$ws.Range("E31").Value = "  +2.01%  "

$ws.Range("D32").Value = "'168.90"
$ws.Range("E32").Value = "  -2.27%  "

$ws.Range("D33").Value = "'20.28"
$ws.Range("E33").Value = "  -1.14%  "

$ws.Range("D34").Value = "'4.92"
$ws.Range("E34").Value = "  +2.82%  "

$ws.Range("E35").Value = "  +0.20%  "

$ws.Range("D36").Value = "'1.45"
$ws.Range("E36").Value = "  +1.05%  "

$ws.Range("D37").Value = "'1.79"
$ws.Range("E37").Value = "  -1.10%  "

$ws.Range("D38").Value = "'0.994"
$ws.Range("E38").Value = "  -0.17%  "

$ws.Range("D39").Value = "'6.14"
$ws.Range("E39").Value = "  +11.18%  "

$ws.Range("D40").Value = "'4.16"
$ws.Range("E40").Value = "  -1.24%  "

$ws.Range("D41").Value = "'329.87"
$ws.Range("E41").Value = "  -4.58%  "

$ws.Range("D42").Value = "'39.43"
$ws.Range("E42").Value = "  +0.96%  "

$ws.Range("D43").Value = "'21.71"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("D44").Value = "'0.0591"
$ws.Range("E44").Value = "  +0.64%  "

$ws.Range("D45").Value = "'21.66"
$ws.Range("E45").Value = "  -0.93%  "

$ws.Range("D46").Value = "'0.634"
$ws.Range("E46").Value = "  -2.06%  "

$ws.Range("D47").Value = "'0.0255"
$ws.Range("E47").Value = "  -0.04%  "

$ws.Range("D48").Value = "'135.87"
$ws.Range("E48").Value = "  -2.31%  "

$ws.Range("D49").Value = "'0.101"
$ws.Range("E49").Value = "  +0.49%  "

$ws.Range("D50").Value = "'0.999"
$ws.Range("E50").Value = "  -0.09%  "

$ws.Range("E51").Value = "  +0.67%  "
